$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (COVID-19 Indicator): mark as selected by lasso method
$ws.Range("B18").Value = "lasso"
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = "lasso"
$ws.Range("G18").Value = 1
